function Set-ParaXml($para, $innerXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

function Insert-EmptyParaAfter($para) {
    $ins = $para.Range.InsertParagraphAfter()
    $newPara = $para.Next()
    Set-ParaXml $newPara '<w:body><w:p/></w:body>'
}

function Insert-ParaXmlBefore($para, $innerXml) {
    $ins = $para.Range.InsertParagraphBefore()
    # After InsertParagraphBefore, $para itself now refers to the freshly
    # created (blank) paragraph; the original content shifted to $para.Next().
    Set-ParaXml $para $innerXml
}

$d = $word.ActiveDocument

# ============================================================
# Work from the END of the document toward the start so that
# paragraph indices of not-yet-processed paragraphs stay valid.
# ============================================================

# --- 12. Delete the entire last "Fiche - Lecture et ecriture de
#     fichiers (NumPy et Pandas)" section (paragraphs 38..46) ---
$total = $d.Paragraphs.Count
$pStart = $d.Paragraphs(38)
$pEnd = $d.Paragraphs($total)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# --- 11.10.9 merge "NumPy"/"numpy" runs that had proofErr wrapping ---

# Paragraph 30: "Operations sur un tableau numpy contenant ou non des valeurs manquantes"
$p30 = $d.Paragraphs(30)
Set-ParaXml $p30 '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Opérations sur un tableau numpy contenant ou non des valeurs manquantes</w:t></w:r></w:p></w:body>'

# Paragraph 29: "Creation de tableaux numpy"
$p29 = $d.Paragraphs(29)
Set-ParaXml $p29 '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Création de tableaux numpy</w:t></w:r></w:p></w:body>'

# Paragraph 28: "Fiche - NumPy, graphiques (barres d'erreur, droite de regression lineaire)"
$p28 = $d.Paragraphs(28)
Set-ParaXml $p28 '<w:body><w:p><w:pPr><w:pStyle w:val="Titre1"/></w:pPr><w:r><w:t>Fiche – NumPy, graphiques (barres d''erreur, droite de régression linéaire)</w:t></w:r></w:p></w:body>'

# --- 8. Remove <w:lastRenderedPageBreak/> from "Modification d'un dictionnaire" ---
$p24 = $d.Paragraphs(24)
Set-ParaXml $p24 '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t>Modification d’un dictionnaire</w:t></w:r></w:p></w:body>'

# --- 7. Insert manual page break paragraph before "Fiche - Dictionnaires" (para 18)
#     and add <w:lastRenderedPageBreak/> marker before its text run ---
$p18 = $d.Paragraphs(18)
Insert-ParaXmlBefore $p18 '<w:body><w:p><w:r><w:br w:type="page"/></w:r></w:p></w:body>'
$p18 = $d.Paragraphs(19)
Set-ParaXml $p18 '<w:body><w:p><w:pPr><w:pStyle w:val="Titre1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Fiche – Dictionnaires</w:t></w:r></w:p></w:body>'

# --- 6. Merge "3. Visualisation de donnees avec Matplotlib" runs ---
$p13 = $d.Paragraphs(13)
Set-ParaXml $p13 '<w:body><w:p><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:r><w:t>3. Visualisation de données avec Matplotlib</w:t></w:r></w:p></w:body>'

# --- 5. "Parcourir" gains a run " une chaine", then an empty paragraph follows ---
$p12 = $d.Paragraphs(12)
Set-ParaXml $p12 '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t>Parcourir</w:t></w:r><w:r><w:t xml:space="preserve"> une chaine</w:t></w:r></w:p></w:body>'
$p12 = $d.Paragraphs(12)
Insert-EmptyParaAfter $p12

# --- 4. "Remplacement" gains a run " d'une chaine" ---
$p11 = $d.Paragraphs(11)
Set-ParaXml $p11 '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t>Remplacement</w:t></w:r><w:r><w:t xml:space="preserve"> d''une chaine</w:t></w:r></w:p></w:body>'

# --- 3. "Conversion" gains a run " d'une chaine" ---
$p10 = $d.Paragraphs(10)
Set-ParaXml $p10 '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t>Conversion</w:t></w:r><w:r><w:t xml:space="preserve"> d''une chaine</w:t></w:r></w:p></w:body>'

# --- 2. "longueur" splits into "L" + "ongueur" + " d'une chaine" ---
$p9 = $d.Paragraphs(9)
Set-ParaXml $p9 '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t>L</w:t></w:r><w:r><w:t>ongueur</w:t></w:r><w:r><w:t xml:space="preserve"> d''une chaine</w:t></w:r></w:p></w:body>'

# --- 1. Insert empty paragraph after "Tris" (paragraph 7) ---
$p7 = $d.Paragraphs(7)
Insert-EmptyParaAfter $p7

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
